$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 48.9
$ws.Range("I11").Value = 48.9
$ws.Range("K11").Value = 48.9
$ws.Range("M11").Value = 91.09999999999999
$ws.Range("H40").Value = 1083.75
$ws.Range("I40").Value = 756.6667
$ws.Range("K40").Value = 756.6667
$ws.Range("M40").Value = -581.6667
$ws.Range("H51").Value = 6599.6
$ws.Range("I51").Value = 8999.333000000001
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 8999.333000000001
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -8515.333000000001
$ws.Range("N51").Value = -3968
$ws.Range("H88").Value = 1178.5555
$ws.Range("I88").Value = 501.5
$ws.Range("J88").Value = 1372
$ws.Range("K88").Value = 501.5
$ws.Range("L88").Value = 1372
$ws.Range("M88").Value = -95.5
$ws.Range("N88").Value = -2184
$ws.Range("H91").Value = 1178.5555
$ws.Range("I91").Value = 501.5
$ws.Range("J91").Value = 1372
$ws.Range("K91").Value = 501.5
$ws.Range("L91").Value = 1372
$ws.Range("M91").Value = 902.5
$ws.Range("N91").Value = -4180
$ws.Range("H106").Value = 2880.6177
$ws.Range("J106").Value = 3486.3809
$ws.Range("L106").Value = 3486.3809
$ws.Range("N106").Value = -4748.3809
$ws.Range("H125").Value = 426.63635
$ws.Range("I125").Value = 366.2857
$ws.Range("J125").Value = 454.8
$ws.Range("K125").Value = 3296.5713
$ws.Range("L125").Value = 4093.2
$ws.Range("M125").Value = -836.5713000000001
$ws.Range("N125").Value = -9013.200000000001
$ws.Range("H129").Value = 271226.25
$ws.Range("J129").Value = 286696.6
$ws.Range("L129").Value = 860089.7999999999
$ws.Range("N129").Value = -870089.7999999999
$ws.Range("H132").Value = 18869820
$ws.Range("I132").Value = 20001918
$ws.Range("J132").Value = 1502
$ws.Range("K132").Value = 60005754
$ws.Range("L132").Value = 4506
$ws.Range("M132").Value = -60003224
$ws.Range("N132").Value = -9566
$ws.Range("H138").Value = 3144.3545
$ws.Range("I138").Value = 3141.5715
$ws.Range("J138").Value = 3144.9539
$ws.Range("K138").Value = 9424.7145
$ws.Range("L138").Value = 9434.861699999999
$ws.Range("M138").Value = -4284.7145
$ws.Range("N138").Value = -19714.8617

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1505.3125
$ws.Range("I2").Value = 1548.9
$ws.Range("J2").Value = 1432.6666
$ws.Range("K2").Value = 1548.9
$ws.Range("L2").Value = 1432.6666
$ws.Range("M2").Value = -1435.9
$ws.Range("N2").Value = -1658.6666
$ws.Range("H32").Value = 10683.906
$ws.Range("I32").Value = 8133.4814
$ws.Range("J32").Value = 24456.2
$ws.Range("K32").Value = 8133.4814
$ws.Range("L32").Value = 24456.2
$ws.Range("M32").Value = -7846.4814
$ws.Range("N32").Value = -25030.2
$ws.Range("H61").Value = 12350383
$ws.Range("I61").Value = 19612818
$ws.Range("J61").Value = 4242.7
$ws.Range("K61").Value = 19612818
$ws.Range("L61").Value = 4242.7
$ws.Range("M61").Value = -19612606
$ws.Range("N61").Value = -4666.7
$ws.Range("H92").Value = 32000
$ws.Range("J92").Value = 32000
$ws.Range("L92").Value = 32000
$ws.Range("N92").Value = -36992
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H116").Value = 1505.3125
$ws.Range("I116").Value = 1548.9
$ws.Range("J116").Value = 1432.6666
$ws.Range("K116").Value = 1548.9
$ws.Range("L116").Value = 1432.6666
$ws.Range("M116").Value = 745.0999999999999
$ws.Range("N116").Value = -6020.6666
$ws.Range("H135").Value = 57768.75
$ws.Range("J135").Value = 57768.75
$ws.Range("L135").Value = 57768.75
$ws.Range("N135").Value = -67908.75
$ws.Range("H136").Value = 12350383
$ws.Range("I136").Value = 19612818
$ws.Range("J136").Value = 4242.7
$ws.Range("K136").Value = 58838454
$ws.Range("L136").Value = 12728.1
$ws.Range("M136").Value = -58835904
$ws.Range("N136").Value = -17828.1

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1505.3125
$ws.Range("I3").Value = 1548.9
$ws.Range("J3").Value = 1432.6666
$ws.Range("K3").Value = 1548.9
$ws.Range("L3").Value = 1432.6666
$ws.Range("M3").Value = -1434.9
$ws.Range("N3").Value = -1660.6666
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H134").Value = 5805.9614
$ws.Range("I134").Value = 5877.2104
$ws.Range("J134").Value = 5612.5713
$ws.Range("K134").Value = 17631.6312
$ws.Range("L134").Value = 16837.7139
$ws.Range("M134").Value = -15096.6312
$ws.Range("N134").Value = -21907.7139

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 17344.625
$ws.Range("I58").Value = 1961.2
$ws.Range("J58").Value = 24337.092
$ws.Range("K58").Value = 1961.2
$ws.Range("L58").Value = 24337.092
$ws.Range("M58").Value = -1758.2
$ws.Range("N58").Value = -24743.092
$ws.Range("H105").Value = 1594.1666
$ws.Range("I105").Value = 1263.625
$ws.Range("K105").Value = 1263.625
$ws.Range("M105").Value = 483.375
$ws.Range("H134").Value = 47620004
$ws.Range("I134").Value = 52632504
$ws.Range("K134").Value = 157897512
$ws.Range("M134").Value = -157894977
$ws.Range("H136").Value = 17344.625
$ws.Range("I136").Value = 1961.2
$ws.Range("J136").Value = 24337.092
$ws.Range("K136").Value = 5883.6
$ws.Range("L136").Value = 73011.276
$ws.Range("M136").Value = -3333.6
$ws.Range("N136").Value = -78111.276

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1174.4546
$ws.Range("I117").Value = 983.8570999999999
$ws.Range("J117").Value = 1508
$ws.Range("K117").Value = 2951.5713
$ws.Range("L117").Value = 4524
$ws.Range("M117").Value = 490.4287000000004
$ws.Range("N117").Value = -11408
$ws.Range("H131").Value = 727.05
$ws.Range("J131").Value = 758.6264
$ws.Range("L131").Value = 2275.8792
$ws.Range("N131").Value = -12355.8792

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4043.75
$ws.Range("I132").Value = 3490.3845
$ws.Range("J132").Value = 5071.4287
$ws.Range("K132").Value = 10471.1535
$ws.Range("L132").Value = 15214.2861
$ws.Range("M132").Value = -7941.1535
$ws.Range("N132").Value = -20274.2861
$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -34860
$ws.Range("H136").Value = 3095.7827
$ws.Range("I136").Value = 3164.5715
$ws.Range("J136").Value = 2988.7778
$ws.Range("K136").Value = 9493.7145
$ws.Range("L136").Value = 8966.3334
$ws.Range("M136").Value = -6943.7145
$ws.Range("N136").Value = -14066.3334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 15000000
$ws.Range("J24").Value = 10000000
$ws.Range("L24").Value = 10000000
$ws.Range("N24").Value = -10000460
$ws.Range("H33").Value = 9800
$ws.Range("J33").Value = 9800
$ws.Range("L33").Value = 9800
$ws.Range("N33").Value = -10300
$ws.Range("H36").Value = 9800
$ws.Range("J36").Value = 9800
$ws.Range("L36").Value = 9800
$ws.Range("N36").Value = -10300
$ws.Range("H40").Value = 9800
$ws.Range("J40").Value = 9800
$ws.Range("L40").Value = 9800
$ws.Range("N40").Value = -10098
$ws.Range("H132").Value = 10001059
$ws.Range("I132").Value = 11905496
$ws.Range("J132").Value = 2762.25
$ws.Range("K132").Value = 35716488
$ws.Range("L132").Value = 8286.75
$ws.Range("M132").Value = -35713958
$ws.Range("N132").Value = -13346.75

Write-Output "Applied all Typhon_Profits updates"
